# "added support for waiting jquery to load"
#
# Update the test-case rows on Sheet1:
#  - row 4 (E4) now calls the new "funcC" test function instead of "funcA"
#  - row 6 (E6) now calls "funcA" instead of the removed "funcB" function
#  - row 5's leftover EOF marker in column B is cleared
#  - the active selection moves to E8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4").Value = "funcC"
$ws.Range("E6").Value = "funcA"
$ws.Range("B5").Value = ""

$ws.Range("E8").Select() | Out-Null
